$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-08 Thursday" "2026-01-09 Friday"

Replace-Text "302÷5=" "730÷6="
Replace-Text "207÷9=" "890÷6="
Replace-Text "864÷5=" "690÷6="
Replace-Text "144÷2=" "270÷3="
Replace-Text "640÷2=" "513÷9="

Replace-Text "303÷8=" "133÷9="
Replace-Text "863÷7=" "693÷9="
Replace-Text "397÷4=" "495÷5="
Replace-Text "455÷6=" "572÷2="
Replace-Text "743÷8=" "149÷6="

Replace-Text "683÷7=" "132÷8="
Replace-Text "106÷7=" "985÷3="
Replace-Text "940÷7=" "389÷7="
Replace-Text "957÷5=" "169÷8="
Replace-Text "799÷2=" "755÷2="

Replace-Text "266÷9=" "198÷4="
Replace-Text "994÷2=" "155÷2="
Replace-Text "354÷5=" "959÷3="
Replace-Text "525÷6=" "182÷5="
Replace-Text "587÷3=" "900÷7="

Replace-Text "414÷9=" "321÷2="
Replace-Text "642÷4=" "178÷6="
Replace-Text "440÷4=" "609÷5="
Replace-Text "694÷8=" "923÷5="
Replace-Text "810÷3=" "683÷8="
